# before.xlsx has a single worksheet named "Magh" with a small logo
# picture anchored at the top-left and the selection sitting on I19.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab "Magh" -> "Sheet"
$ws.Name = "Sheet"

# 2. The logo picture (Shapes(1)) shrinks by a hair: its extent goes
#    from cx=1245960 / cy=373320 EMU to cx=1245600 / cy=372960 EMU
#    while its top-left anchor (col0/99360, row0/79200) stays put.
#    1 point = 12700 EMU, so:
#      1245600 EMU -> 98.07874015748031 pt
#       372960 EMU -> 29.366929133858267 pt
if ($ws.Shapes.Count -ge 1) {
    $shp = $ws.Shapes.Item(1)
    $shp.Width = 98.07874015748031
    $shp.Height = 29.366929133858267
}

# 3. Move the active cell / selection from I19 to W22
$ws.Range("W22").Select()
